# Insert one new data row (week) before the existing row 718 on the
# "Mandarina" sheet, shifting the remaining rows (718..796) down to
# (719..797), and populate the new row with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 718:796 down by one row to make room for the new record.
$ws.Rows("718:718").Insert()

# Fill in the newly inserted row 718 with the new observation.
$ws.Range("A718").Value2 = 10
$ws.Range("B718").Value2 = "Vega Modelo de Temuco"
$ws.Range("C718").Value2 = "La Araucanía"
$ws.Range("D718").Value2 = 44951
$ws.Range("E718").Value2 = 9
$ws.Range("F718").Value2 = "Fruta"
$ws.Range("G718").Value2 = 100102
$ws.Range("H718").Value2 = "Cítricos"
$ws.Range("I718").Value2 = 100102004
$ws.Range("J718").Value2 = "Mandarina"
$ws.Range("K718").Value2 = "Murcott"
$ws.Range("L718").Value2 = "Primera"
$ws.Range("M718").Value2 = 100
$ws.Range("N718").Value2 = 12000
$ws.Range("O718").Value2 = 12000
$ws.Range("P718").Value2 = 12000
$ws.Range("Q718").Value2 = "$/bandeja 10 kilos"
$ws.Range("R718").Value2 = "Región de O'Higgins"
$ws.Range("S718").Value2 = 1200
$ws.Range("T718").Value2 = 10
